$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Paragraph 1: "Yerba Buena, 29 de Octubre de 1992"
#   + keepNext, spacing after 200 -> 240
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Format.KeepWithNext = 1
$p1.Format.SpaceAfter = 12

# ---------------------------------------------------------------------------
# Paragraph 2: "ORDENANZA Nº 505"
#   + keepNext, spacing before 240 after 360, bold
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$p2.Format.KeepWithNext = 1
$p2.Format.SpaceBefore = 12
$p2.Format.SpaceAfter = 18
$p2.Range.Font.Bold = 1

# ---------------------------------------------------------------------------
# Paragraph 3: "VISTO: Que entre el Personal Municipal..."
#   Split into two paragraphs:
#     3a: "VISTO: "  (bold, keepNext, before 240 after 120)
#     3b: " " + "Que entre el Personal..." (keepNext, after 120)
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$p3start = $p3.Range.Start
$splitPos = $p3start + 7  # length of "VISTO: "
$splitRange = $d.Range($splitPos, $splitPos)
$splitRange.InsertParagraphAfter()

# New paragraph (continuation) now begins right after the inserted mark.
$p3b = $d.Paragraphs.Item(4)
$p3bStart = $p3b.Range.Start
$insertPoint = $d.Range($p3bStart, $p3bStart)
$insertPoint.InsertAfter(" ")
# Force the freshly inserted space to stay a distinct run by toggling a
# character attribute on/off (otherwise it silently merges into the
# following run because the formatting ends up identical).
$spaceRange = $d.Range($p3bStart, $p3bStart + 1)
$spaceRange.Font.Bold = 1
$spaceRange.Font.Bold = 0

$p3a = $d.Paragraphs.Item(3)
$p3a.Format.KeepWithNext = 1
$p3a.Format.SpaceBefore = 12
$p3a.Format.SpaceAfter = 6
$p3a.Range.Font.Bold = 1

$p3b = $d.Paragraphs.Item(4)
$p3b.Format.KeepWithNext = 1
$p3b.Format.SpaceAfter = 6

# ---------------------------------------------------------------------------
# Paragraphs that follow (originally 4..7, now shifted by +1 to 5..8):
#   "Que se hace necesario..."  / "CONIDERANDO:" /
#   "Que esta norma..." / "Que de esta forma..."
#   + keepNext, spacing after 200 -> 120
# ---------------------------------------------------------------------------
foreach ($idx in 5,6,7,8) {
    $p = $d.Paragraphs.Item($idx)
    $p.Format.KeepWithNext = 1
    $p.Format.SpaceAfter = 6
}

# ---------------------------------------------------------------------------
# Paragraph (originally 8, now 9): "POR EL CONCEJO DELIBERANTE SANCIONA..."
#   + keepNext, spacing before 360 after 360, indent left/right 1984,
#   bold, and the two runs ("POR " / "EL CONCEJO...") merge into one run.
# ---------------------------------------------------------------------------
$searchText = "POR EL CONCEJO DELIBERANTE SANCIONA CON FUERZA DE ORDENANZA"
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, $searchText, 2) | Out-Null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd() -eq $searchText) {
        $p.Format.KeepWithNext = 1
        $p.Format.SpaceBefore = 18
        $p.Format.SpaceAfter = 18
        $p.Format.LeftIndent = 99.2
        $p.Format.RightIndent = 99.2
        $p.Range.Font.Bold = 1
    }
}

# ---------------------------------------------------------------------------
# ARTICULO PRIMERO / SEGUNDO / TERCERO / CUARTO / QUINTO paragraphs:
#   + keepNext, spacing after 200 -> 120
#   "ARTICULO X: " -> "ARTICULO X:" (underlined) + " " (plain) as 2 runs
# ---------------------------------------------------------------------------
$labels = "ARTICULO PRIMERO", "ARTICULO SEGUNDO", "ARTICULO TERCERO", "ARTICULO CUARTO", "ARTICULO QUINTO"
foreach ($label in $labels) {
    $full = $label + ": "
    $find2 = $d.Content.Find
    $find2.ClearFormatting()
    $find2.Text = $full
    $find2.Execute() | Out-Null
    $r = $find2.Parent
    $rWord = $d.Range($r.Start, $r.End - 1)
    $rWord.Font.Underline = 1

    # NOTE: a Paragraph fetched off a narrow sub-Range's .Paragraphs
    # collection is a detached wrapper whose Format writes don't make it
    # back into the real document; walk $d.Paragraphs directly instead.
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Start -le $r.Start -and $p.Range.End -ge $r.End) {
            $p.Format.KeepWithNext = 1
            $p.Format.SpaceAfter = 6
        }
    }
}
